# Applies the "maleState" column edits to Sheet1.
#   - Adds a new header cell G1 = "maleState"
#   - Replaces the old G2 legend ("lead/fol") with the coded legend text
#   - Converts the various G-column "L"/"F" string flags into numeric codes
#     (lea = 3, fol = 2, sol = 1, juvsol = 0) and fills in previously blank
#     G cells with the appropriate numeric code
#   - Grows row 2's height to fit the new wrapped legend text
#   - Leaves the last selected cell as G67 (where the user ended up)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header / legend text
$ws.Range("G1").Value = "maleState"
$ws.Range("G2").Value = "lea = 3, fol = 2, sol = 1, juvsol = 0"

# Row 2 needs to grow to fit the two-line wrapped legend text
$ws.Rows.Item(2).RowHeight = 58.2

# Existing "L"/"F" string codes -> numeric codes
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 2
$ws.Range("G19").Value = 3
$ws.Range("G20").Value = 2
$ws.Range("G60").Value = 3
$ws.Range("G61").Value = 2

# Previously empty G cells now carry numeric codes
$ws.Range("G6").Value = 0
$ws.Range("G13").Value = 3
$ws.Range("G25").Value = 0
$ws.Range("G32").Value = 3
$ws.Range("G38").Value = 3
$ws.Range("G42").Value = 1
$ws.Range("G43").Value = 3
$ws.Range("G45").Value = 0
$ws.Range("G52").Value = 3
$ws.Range("G67").Value = 1

# Leave the selection where the editor last left it
$ws.Activate()
$ws.Range("G67").Select()
